$wb = $excel.ActiveWorkbook

# Sheet ALC, row 70 (Leve Item ID 12604)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 2004.091
$ws.Cells.Item(70, 9).Value = 2072.7334
$ws.Cells.Item(70, 10).Value = 1857
$ws.Cells.Item(70, 11).Value = 6218.2002
$ws.Cells.Item(70, 12).Value = 5571
$ws.Cells.Item(70, 13).Value = -5948.2002
$ws.Cells.Item(70, 14).Value = -6111

# Sheet ALC, row 73 (Leve Item ID 12604)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 2004.091
$ws.Cells.Item(73, 9).Value = 2072.7334
$ws.Cells.Item(73, 10).Value = 1857
$ws.Cells.Item(73, 11).Value = 6218.2002
$ws.Cells.Item(73, 12).Value = 5571
$ws.Cells.Item(73, 13).Value = -5282.2002
$ws.Cells.Item(73, 14).Value = -7443

# Sheet ALC, row 98 (Leve Item ID 36237)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 577.3461
$ws.Cells.Item(98, 9).Value = 577.3461
$ws.Cells.Item(98, 11).Value = 577.3461
$ws.Cells.Item(98, 13).Value = 920.6539

# Sheet ALC, row 122 (Leve Item ID 36237)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 577.3461
$ws.Cells.Item(122, 9).Value = 577.3461
$ws.Cells.Item(122, 11).Value = 1732.0383
$ws.Cells.Item(122, 13).Value = 717.9617000000001

# Sheet ALC, row 134 (Leve Item ID 41997)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(134, 8).Value = 139105.2
$ws.Cells.Item(134, 10).Value = 139105.2
$ws.Cells.Item(134, 12).Value = 139105.2
$ws.Cells.Item(134, 14).Value = -149245.2

# Sheet ALC, row 137 (Leve Item ID 44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 5285.0312
$ws.Cells.Item(137, 9).Value = 19346
$ws.Cells.Item(137, 11).Value = 58038
$ws.Cells.Item(137, 13).Value = -55488

# Sheet ARM, row 2 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 35240624
$ws.Cells.Item(2, 9).Value = 52858784
$ws.Cells.Item(2, 11).Value = 52858784
$ws.Cells.Item(2, 13).Value = -52858671

# Sheet ARM, row 92 (Leve Item ID 18050)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(92, 8).Value = 159999.5
$ws.Cells.Item(92, 9).Value = 120000
$ws.Cells.Item(92, 10).Value = 199999
$ws.Cells.Item(92, 11).Value = 120000
$ws.Cells.Item(92, 12).Value = 199999
$ws.Cells.Item(92, 13).Value = -117504
$ws.Cells.Item(92, 14).Value = -204991

# Sheet ARM, row 116 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 35240624
$ws.Cells.Item(116, 9).Value = 52858784
$ws.Cells.Item(116, 11).Value = 52858784
$ws.Cells.Item(116, 13).Value = -52856490

# Sheet ARM, row 119 (Leve Item ID 26287)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(119, 8).Value = 71428.28999999999
$ws.Cells.Item(119, 10).Value = 71428.28999999999
$ws.Cells.Item(119, 12).Value = 71428.28999999999
$ws.Cells.Item(119, 14).Value = -81104.28999999999

# Sheet ARM, row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1734.7872
$ws.Cells.Item(132, 9).Value = 1573.8055
$ws.Cells.Item(132, 10).Value = 2261.6365
$ws.Cells.Item(132, 11).Value = 4721.416499999999
$ws.Cells.Item(132, 12).Value = 6784.9095
$ws.Cells.Item(132, 13).Value = -2191.416499999999
$ws.Cells.Item(132, 14).Value = -11844.9095

# Sheet BSM, row 3 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 35240624
$ws.Cells.Item(3, 9).Value = 52858784
$ws.Cells.Item(3, 11).Value = 52858784
$ws.Cells.Item(3, 13).Value = -52858670

# Sheet BSM, row 20 (Leve Item ID 14149)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2708.743
$ws.Cells.Item(20, 9).Value = 2345.926
$ws.Cells.Item(20, 11).Value = 2345.926
$ws.Cells.Item(20, 13).Value = -2098.926

# Sheet BSM, row 86 (Leve Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2902
$ws.Cells.Item(86, 9).Value = 2863.4
$ws.Cells.Item(86, 11).Value = 2863.4
$ws.Cells.Item(86, 13).Value = -1740.4

# Sheet BSM, row 89 (Leve Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 2902
$ws.Cells.Item(89, 9).Value = 2863.4
$ws.Cells.Item(89, 11).Value = 14317
$ws.Cells.Item(89, 13).Value = -8701

# Sheet CRP, row 99 (Leve Item ID 36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 9237.119000000001
$ws.Cells.Item(99, 10).Value = 14484.263
$ws.Cells.Item(99, 12).Value = 14484.263
$ws.Cells.Item(99, 14).Value = -17480.263

# Sheet CRP, row 105 (Leve Item ID 19928)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 3388.7058
$ws.Cells.Item(105, 9).Value = 1576.7858
$ws.Cells.Item(105, 11).Value = 1576.7858
$ws.Cells.Item(105, 13).Value = 170.2141999999999

# Sheet CRP, row 124 (Leve Item ID 34285)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(124, 8).Value = 50401.668
$ws.Cells.Item(124, 10).Value = 50401.668
$ws.Cells.Item(124, 12).Value = 50401.668
$ws.Cells.Item(124, 14).Value = -55311.668

# Sheet CRP, row 126 (Leve Item ID 36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 9237.119000000001
$ws.Cells.Item(126, 10).Value = 14484.263
$ws.Cells.Item(126, 12).Value = 43452.789
$ws.Cells.Item(126, 14).Value = -48392.789

# Sheet CRP, row 132 (Leve Item ID 44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 5078.1665
$ws.Cells.Item(132, 9).Value = 2131.5789
$ws.Cells.Item(132, 11).Value = 6394.736699999999
$ws.Cells.Item(132, 13).Value = -3864.736699999999

# Sheet CUL, row 92 (Leve Item ID 19841)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 742.5
$ws.Cells.Item(92, 9).Value = 764.75
$ws.Cells.Item(92, 10).Value = 720.25
$ws.Cells.Item(92, 11).Value = 2294.25
$ws.Cells.Item(92, 12).Value = 2160.75
$ws.Cells.Item(92, 13).Value = -1046.25
$ws.Cells.Item(92, 14).Value = -4656.75

# Sheet CUL, row 97 (Leve Item ID 19846)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(97, 8).Value = 292.13333
$ws.Cells.Item(97, 9).Value = 160
$ws.Cells.Item(97, 11).Value = 480
$ws.Cells.Item(97, 13).Value = 16

# Sheet CUL, row 139 (Leve Item ID 44102)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 3906.8462
$ws.Cells.Item(139, 9).Value = 1272.1111
$ws.Cells.Item(139, 10).Value = 9835
$ws.Cells.Item(139, 11).Value = 3816.3333
$ws.Cells.Item(139, 12).Value = 29505
$ws.Cells.Item(139, 13).Value = 1323.6667
$ws.Cells.Item(139, 14).Value = -39785

# Sheet GSM, row 43 (Leve Item ID 4218)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 26400
$ws.Cells.Item(43, 9).Value = 20000
$ws.Cells.Item(43, 10).Value = 36000
$ws.Cells.Item(43, 11).Value = 20000
$ws.Cells.Item(43, 12).Value = 36000
$ws.Cells.Item(43, 13).Value = -19849
$ws.Cells.Item(43, 14).Value = -36302

# Sheet GSM, row 102 (Leve Item ID 36169)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 4989.1816
$ws.Cells.Item(102, 9).Value = 5152.923
$ws.Cells.Item(102, 10).Value = 4752.6665
$ws.Cells.Item(102, 11).Value = 5152.923
$ws.Cells.Item(102, 12).Value = 4752.6665
$ws.Cells.Item(102, 13).Value = -3530.923
$ws.Cells.Item(102, 14).Value = -7996.6665

# Sheet GSM, row 122 (Leve Item ID 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 8165.75
$ws.Cells.Item(122, 9).Value = 7651.3076
$ws.Cells.Item(122, 10).Value = 9121.143
$ws.Cells.Item(122, 11).Value = 22953.9228
$ws.Cells.Item(122, 12).Value = 27363.429
$ws.Cells.Item(122, 13).Value = -20503.9228
$ws.Cells.Item(122, 14).Value = -32263.429

# Sheet GSM, row 132 (Leve Item ID 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 5063.091
$ws.Cells.Item(132, 9).Value = 4740.737
$ws.Cells.Item(132, 11).Value = 14222.211
$ws.Cells.Item(132, 13).Value = -11692.211

# Sheet GSM, row 138 (Leve Item ID 42325)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(138, 8).Value = 138943.2
$ws.Cells.Item(138, 10).Value = 138943.2
$ws.Cells.Item(138, 12).Value = 138943.2
$ws.Cells.Item(138, 14).Value = -149223.2

# Sheet LTW, row 55 (Leve Item ID 5284)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 677.8570999999999
$ws.Cells.Item(55, 9).Value = 389.81818
$ws.Cells.Item(55, 11).Value = 389.81818
$ws.Cells.Item(55, 13).Value = -216.81818

# Sheet LTW, row 136 (Leve Item ID 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 4413.6
$ws.Cells.Item(136, 9).Value = 2244.5789
$ws.Cells.Item(136, 11).Value = 6733.736699999999
$ws.Cells.Item(136, 13).Value = -4183.736699999999

# Sheet WVR, row 124 (Leve Item ID 34280)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(124, 14).ClearContents()

# Sheet WVR, row 127 (Leve Item ID 35414)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(127, 8).Value = 0
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 14).ClearContents()

# Sheet WVR, row 132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2227.6304
$ws.Cells.Item(132, 9).Value = 2177.0698
$ws.Cells.Item(132, 11).Value = 6531.209400000001
$ws.Cells.Item(132, 13).Value = -4001.209400000001

# Sheet WVR, row 135 (Leve Item ID 42043)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(135, 8).Value = 250285.25
$ws.Cells.Item(135, 10).Value = 250285.25
$ws.Cells.Item(135, 12).Value = 250285.25
$ws.Cells.Item(135, 14).Value = -260425.25
